$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$t = $m.Theme
try {
    $t.Name = "Integral"
    Write-Host "set ok"
} catch {
    Write-Host "ERR: $_"
}
